$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels to the new, clearer names.
$ws.Range("F1").Value = "repo"
$ws.Range("G1").Value = "languages"
$ws.Range("H1").Value = "libraries_tools"

# Replace the "concepts" metadata for the Credit Risk row (row 8) -- swap in the
# new page's metadata content.
$ws.Range("I8").Value = "linear regression, t-test, hypothesis testing, study design"

# Row 8 no longer needs the taller wrapped height; let it size back to default.
$ws.Rows.Item(8).EntireRow.AutoFit()

# Scroll the view so column F is the left-most visible column, and select I8
# (the "expand language section button" target cell) as the single active cell.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("I8").Select()
